# Generate Report for Handoff
#
# c0a03d35-d1ce-4013-b7e1-168c9054c3ea just reached "Ready for handoff" status
# with a fresh handoff timestamp. The localization-status report re-sorts its
# active rows (6-9) accordingly: the three still-in-progress / already-ready
# rows shift up one position, and c0a03d35's row is pushed to the bottom of
# the block with its refreshed status/date.

$wb = $excel.ActiveWorkbook

function Rotate-Sheet($SheetName, $Columns, $StatusColumns, $DateColumn, $NewDateValue, $NewStatusValue) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Snapshot current values for rows 6-9 across the columns we care about.
    $orig = @{}
    for ($r = 6; $r -le 9; $r++) {
        $rowVals = @{}
        foreach ($col in $Columns) {
            $rowVals[$col] = $ws.Range("$col$r").Value2
        }
        $orig[$r] = $rowVals
    }

    # Rows 6,7,8 take on the content that used to live one row below them.
    for ($r = 6; $r -le 8; $r++) {
        $below = $orig[$r + 1]
        foreach ($col in $Columns) {
            $ws.Range("$col$r").Value = $below[$col]
        }
    }

    # Row 9 takes the content that used to be in row 6 (the file that just
    # completed handoff), except for the Status column(s) and Date column,
    # which get refreshed values.
    $top = $orig[6]
    foreach ($col in $Columns) {
        $ws.Range("${col}9").Value = $top[$col]
    }
    foreach ($col in $StatusColumns) {
        $ws.Range("${col}9").Value = $NewStatusValue
    }
    $dateCell = $DateColumn + "9"
    $ws.Range($dateCell).Value = $NewDateValue

    # Update hyperlink display text so it follows the cell content that now
    # sits in each row (the hyperlink targets/r:id stay pinned to their row).
    foreach ($hl in $ws.Hyperlinks) {
        $row = $hl.Range.Row
        $col = $hl.Range.Column
        if ($row -ge 6 -and $row -le 9) {
            $addr = $ws.Cells.Item($row, $col).Address($false, $false)
            $colLetter = ($addr -replace '[0-9]', '')
            if ($Columns -contains $colLetter) {
                $cellAddr = $colLetter + $row
                $hl.TextToDisplay = $ws.Range($cellAddr).Value2
            }
        }
    }
}

# Overview sheet: File Name (A), zh-cn Status (B), de-de Status (C), Latest Handoff Date (D)
Rotate-Sheet "Overview" @("A","B","C","D") @("B","C") "D" "2016-03-22 12:48:57" "Ready for handoff"

# zh-cn sheet: Source File Name (A), File Extension (B), Status (C), Latest Handoff File (D),
# Latest Handoff Datetime (E), Latest Handback DateTime (H), Handoff Reason (J)
Rotate-Sheet "zh-cn" @("A","B","C","D","E","H","J") @("C") "E" "2016-03-22 12:48:53" "Ready for handoff"

# de-de sheet: same shape as zh-cn
Rotate-Sheet "de-de" @("A","B","C","D","E","H","J") @("C") "E" "2016-03-22 12:48:57" "Ready for handoff"
